# "Added 1.1.0 of term"
# Bumps the Version metadata row from 1.0.0 -> 1.1.0 and refreshes the
# Date metadata row to the new publish timestamp, on the "Metadata" sheet.
#
# Sheet layout (column A = property name, column B = value):
#   Row 3: Version | 1.0.0
#   Row 8: Date    | 2023-06-07T11:52:14+02:00

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
